$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B-E
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values for columns B-E
$ws.Range("B2").Value = 55.699168891699358
$ws.Range("C2").Value = 60.778775624858106
$ws.Range("D2").Value = 51.337072624918413
$ws.Range("E2").Value = 71.519478132967635

# Update row 3 values for columns B-E
$ws.Range("B3").Value = 31.908378433147845
$ws.Range("C3").Value = 62.140063266426878
$ws.Range("D3").Value = 54.722123549108382
$ws.Range("E3").Value = 82.065695361940101

# Update the selection to match the new active range
$ws.Range("B1:E3").Select()
